# Update computed power-flow results on Sheet1 for the "380 kV" case.
# Columns C, D, E, F, G, H, J, O for rows 2-25 are replaced with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each inner array: row, C, D, E, F, G, H, J, O
$data = @(
    @(2,0.04969111061828357,0.2107476850597223,0.1680132064088546,1.276671991392575,0.7614461594884006,0.7362247279381791,0.1838998710622519,3.002323868125757),
    @(3,0.04408939871433404,0.2101148881426838,0.1653876853704439,1.235354089528201,0.7201967886405214,0.7216148393470121,0.1786267740461298,2.883864273825338),
    @(4,0.04066648080546997,0.2098108844848454,0.163866114486094,1.210749221148646,0.6953400628318036,0.7130539717894635,0.1755003753982791,2.812967843030606),
    @(5,0.03927570747256937,0.2097082792852802,0.1632687989183026,1.200914363468556,0.6853285548017567,0.7096682597508845,0.1742542540484706,2.784537893400568),
    @(6,0.03904501583093634,0.2096925275957631,0.1631709877360947,1.199292866790941,0.6836732519808351,0.7091122775099876,0.174049020107006,2.779844922308769),
    @(7,0.04064770786612826,0.2098094145292961,0.1638579668473596,1.210615808670951,0.6952045678090997,0.7130078943511364,0.1754834568576626,2.812582562290174),
    @(8,0.04775616062535448,0.2105119633120083,0.1670891265870758,1.262266643986962,0.7471253876040294,0.7311021441241508,0.1820585575032538,2.96109663284426),
    @(9,0.06183100901742478,0.2125597742663388,0.174145148148952,1.369645234141359,0.8527049597045959,0.7698429204376112,0.1958398391220157,3.266996416250322),
    @(10,0.07226051778398812,0.2144724038749786,0.1797709608776827,1.452295328595596,0.9326177770774393,0.8003065429870162,0.2065133899800884,3.50081688383699),
    @(11,0.07702584964360426,0.215431072439884,0.1824269134105236,1.490721737755223,0.9694922927974687,0.8146031334345594,0.2114898347394103,3.60918982072053),
    @(12,0.07883345327934421,0.2158068242919313,0.1834466062456741,1.505392597165354,0.9835314671215656,0.8200801398339195,0.2133917932219589,3.650518376122136),
    @(13,0.07844401556081948,0.2157253337205844,0.1832263766266991,1.50222763760155,0.980504508558397,0.818897755549358,0.2129813932858298,3.641604607028341),
    @(14,0.07717450022919081,0.2154617308370348,0.1825105245772676,1.491926318211341,0.970645785818931,0.8150524626720426,0.2116459590371704,3.612584123659644),
    @(15,0.07639728811501811,0.2153019231383695,0.1820738609860229,1.485632054026922,0.9646168950918366,0.812705344101829,0.2108302469394232,3.594846060815144),
    @(16,0.07194951869286115,0.2144115341675104,0.1795993373710516,1.449800784119191,0.9302184883761697,0.7993810632741543,0.2061906098192026,3.493774967662773),
    @(17,0.06922635870340343,0.2138879912362484,0.1781060963177268,1.42803199014709,0.9092502194446865,0.791319464521905,0.2033753944019594,3.432286183054828),
    @(18,0.06766203167595108,0.2135952039431004,0.1772563275652033,1.415589105775723,0.897238903554836,0.7867239148270926,0.2017675400542061,3.397108262714369),
    @(19,0.06713271233509488,0.2134975042291671,0.1769701731879394,1.4113895324179,0.8931804903970431,0.7851750244724087,0.2012251005844519,3.385230004966104),
    @(20,0.06951604009738332,0.2139428601577862,0.1782641119323713,1.430341244444207,0.9114772440300101,0.7921733633517363,0.2036739001437837,3.438812218089765),
    @(21,0.07754730361043016,0.2155388121710473,0.1827204090742924,1.494948815327547,0.9735394754995355,0.8161802021211031,0.2120377331520729,3.621100258499723),
    @(22,0.08281416027634236,0.2166560221208442,0.1857141314301245,1.537871016488253,1.014541625562515,0.8322384904946318,0.2176059498935814,3.741927523114384),
    @(23,0.08000147273300229,0.2160529654801309,0.1841088788204743,1.514898655921428,0.9926174766380313,0.8236341221844725,0.2146247283204872,3.677284467670063),
    @(24,0.06938507123919635,0.2139180283620732,0.1781926459288243,1.429297005451531,0.9104702703839962,0.791787193684911,0.2035389125313571,3.435861259504236),
    @(25,0.0580082653981151,0.2119340913197121,0.172158948687283,1.339939663534125,0.8237349345114922,0.7590122234520038,0.1920158612539069,3.18265849722826)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value  = $row[1]   # C
    $ws.Cells.Item($r, 4).Value  = $row[2]   # D
    $ws.Cells.Item($r, 5).Value  = $row[3]   # E
    $ws.Cells.Item($r, 6).Value  = $row[4]   # F
    $ws.Cells.Item($r, 7).Value  = $row[5]   # G
    $ws.Cells.Item($r, 8).Value  = $row[6]   # H
    $ws.Cells.Item($r, 10).Value = $row[7]   # J
    $ws.Cells.Item($r, 15).Value = $row[8]   # O
}
